# Applies the "Updated cryptos list" price/volume refresh described in the diff.
# D-column "Price" values are plain numeric-looking strings (e.g. "1.00", "10.60") that
# must stay as literal text (matching the original inlineStr cells), so each such value is
# written with a leading apostrophe (forces text entry) and the cell style is reset back to
# "Normal" afterwards so no numeric formatting/style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    if ($text -match "^[+-]?[0-9]+(\.[0-9]+)?$") {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue "D2" '51.150.60'
Set-TextValue "E2" '  -0.91%  '
Set-TextValue "D3" '3.062.15'
Set-TextValue "E3" '  +1.09%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -0.04%  '
Set-TextValue "D5" '392.02'
Set-TextValue "E5" '  +2.07%  '
Set-TextValue "D6" '101.57'
Set-TextValue "E7" '  -2.18%  '
Set-TextValue "D9" '0.579'
Set-TextValue "E9" '  -1.85%  '
Set-TextValue "D10" '36.71'
Set-TextValue "E10" '  +0.07%  '
Set-TextValue "E11" '  +0.24%  '
Set-TextValue "D12" '0.0847'
Set-TextValue "E12" '  -1.67%  '
Set-TextValue "D13" '3.542.31'
Set-TextValue "E13" '  +1.10%  '
Set-TextValue "D14" '18.25'
Set-TextValue "E14" '  -1.42%  '
Set-TextValue "D15" '7.66'
Set-TextValue "E15" '  -0.91%  '
Set-TextValue "D16" '3.060.11'
Set-TextValue "E16" '  +1.21%  '
Set-TextValue "D17" '0.989'
Set-TextValue "E17" '  +1.73%  '
Set-TextValue "D18" '10.60'
Set-TextValue "E18" '  -0.53%  '
Set-TextValue "D19" '51.127.05'
Set-TextValue "E19" '  -0.96%  '
Set-TextValue "D20" '3.17'
Set-TextValue "E20" '  +3.07%  '
Set-TextValue "E21" '  -0.89%  '
Set-TextValue "D22" '12.22'
Set-TextValue "E22" '  -1.98%  '
Set-TextValue "D23" '69.58'
Set-TextValue "E23" '  -0.53%  '
Set-TextValue "D24" '263.99'
Set-TextValue "E24" '  -1.22%  '
Set-TextValue "D25" '3.13'
Set-TextValue "E25" '  -1.26%  '
Set-TextValue "D26" '7.88'
Set-TextValue "E26" '  -6.63%  '
Set-TextValue "D27" '26.71'
Set-TextValue "E27" '  +1.87%  '
Set-TextValue "D28" '1.00'
Set-TextValue "E28" '  -0.03%  '
Set-TextValue "D29" '7.10'
Set-TextValue "E29" '  -5.77%  '
Set-TextValue "E30" '  -5.62%  '
Set-TextValue "E31" '  -3.15%  '
Set-TextValue "D33" '0.0486'
Set-TextValue "E33" '  +8.37%  '
Set-TextValue "D34" '35.65'
Set-TextValue "E34" '  +4.57%  '
Set-TextValue "E35" '  -0.39%  '
Set-TextValue "D36" '49.95'
Set-TextValue "E36" '  -1.12%  '
Set-TextValue "E37" '  -0.06%  '
Set-TextValue "D38" '3.35'
Set-TextValue "E38" '  +1.41%  '
Set-TextValue "D39" '0.291'
Set-TextValue "E39" '  -2.77%  '
Set-TextValue "D40" '129.56'
Set-TextValue "E40" '  +1.45%  '
Set-TextValue "D41" '16.52'
Set-TextValue "E41" '  -3.11%  '
Set-TextValue "D42" '1.83'
Set-TextValue "E42" '  -1.60%  '
Set-TextValue "E43" '  -1.34%  '
Set-TextValue "B44" 'NEARProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D44" '3.76'
Set-TextValue "E44" '  +1.18%  '
Set-TextValue "B45" 'Stacks'
Set-TextValue "C45" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D45" '2.48'
Set-TextValue "E45" '  -1.84%  '
Set-TextValue "D46" '21.67'
Set-TextValue "E46" '  +0.59%  '
Set-TextValue "E47" '  +3.38%  '
Set-TextValue "E48" '  -0.52%  '
Set-TextValue "D49" '2.062.54'
Set-TextValue "E49" '  +1.86%  '
Set-TextValue "D50" '0.0322'
Set-TextValue "E50" '  +3.04%  '
Set-TextValue "D51" '0.886'
Set-TextValue "E51" '  +11.46%  '
